{"js": "const replacements = [\n  [\"83\u00d788=\", \"75\u00d754=\"],\n  [\"85\u00d793=\", \"97\u00d740=\"],\n  [\"74\u00d788=\", \"57\u00d794=\"],\n  [\"54\u00d718=\", \"36\u00d760=\"],\n  [\"51\u00d785=\", \"30\u00d760=\"],\n  [\"23\u00d732=\", \"66\u00d736=\"],\n  [\"91\u00d725=\", \"94\u00d732=\"],\n  [\"26\u00d728=\", \"93\u00d782=\"],\n  [\"13\u00d721=\", \"97\u00d721=\"],\n  [\"79\u00d714=\", \"64\u00d761=\"],\n  [\"50\u00d785=\", \"92\u00d761=\"],\n  [\"55\u00d799=\", \"42\u00d762=\"],\n  [\"20\u00d788=\", \"21\u00d725=\"],\n  [\"25\u00d760=\", \"51\u00d743=\"],\n  [\"34\u00d755=\", \"79\u00d757=\"],\n  [\"94\u00d733=\", \"47\u00d729=\"],\n  [\"48\u00d760=\", \"36\u00d728=\"],\n  [\"15\u00d774=\", \"22\u00d777=\"],\n  [\"80\u00d789=\", \"61\u00d738=\"],\n  [\"20\u00d716=\", \"62\u00d716=\"],\n  [\"73\u00d741=\", \"13\u00d732=\"],\n  [\"46\u00d782=\", \"67\u00d744=\"],\n  [\"93\u00d764=\", \"11\u00d757=\"],\n  [\"26\u00d722=\", \"11\u00d767=\"],\n  [\"29\u00d795=\", \"27\u00d779=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"83\u00d788=\", \"75\u00d754=\"),\n    @(\"85\u00d793=\", \"97\u00d740=\"),\n    @(\"74\u00d788=\", \"57\u00d794=\"),\n    @(\"54\u00d718=\", \"36\u00d760=\"),\n    @(\"51\u00d785=\", \"30\u00d760=\"),\n    @(\"23\u00d732=\", \"66\u00d736=\"),\n    @(\"91\u00d725=\", \"94\u00d732=\"),\n    @(\"26\u00d728=\", \"93\u00d782=\"),\n    @(\"13\u00d721=\", \"97\u00d721=\"),\n    @(\"79\u00d714=\", \"64\u00d761=\"),\n    @(\"50\u00d785=\", \"92\u00d761=\"),\n    @(\"55\u00d799=\", \"42\u00d762=\"),\n    @(\"20\u00d788=\", \"21\u00d725=\"),\n    @(\"25\u00d760=\", \"51\u00d743=\"),\n    @(\"34\u00d755=\", \"79\u00d757=\"),\n    @(\"94\u00d733=\", \"47\u00d729=\"),\n    @(\"48\u00d760=\", \"36\u00d728=\"),\n    @(\"15\u00d774=\", \"22\u00d777=\"),\n    @(\"80\u00d789=\", \"61\u00d738=\"),\n    @(\"20\u00d716=\", \"62\u00d716=\"),\n    @(\"73\u00d741=\", \"13\u00d732=\"),\n    @(\"46\u00d782=\", \"67\u00d744=\"),\n    @(\"93\u00d764=\", \"11\u00d757=\"),\n    @(\"26\u00d722=\", \"11\u00d767=\"),\n    @(\"29\u00d795=\", \"27\u00d779=\"),\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair[1], $wdReplaceAll) | Out-Null\n}"}
